# Update computed investment-capacity results for the 2025, 2030 and 2035
# sheets with freshly recalculated values received from the server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.002168060809513972
$ws.Range("E2").Value = 0.3759634166263097
$ws.Range("G2").Value = 0.2494892361375063
$ws.Range("I2").Value = 0.3511355
$ws.Range("L2").Value = 0.6154939666666667
$ws.Range("M2").Value = 0.08148166666666669
$ws.Range("N2").Value = 12.70991232200972
$ws.Range("O2").Value = 3.512978867815095

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.04814661662630972
$ws.Range("E2").Value = 0.2278274343708592
$ws.Range("I2").Value = 0.3478968171620428
$ws.Range("L2").Value = 0.1883974735243333
$ws.Range("M2").Value = 0.04652328729864512
$ws.Range("N2").Value = 5.392237274658115
$ws.Range("O2").Value = 2.368284394810271

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.08592412015137646
$ws.Range("B2").Value = 0.02745799999999977
$ws.Range("E2").Value = 0.17352283365371
$ws.Range("I2").Value = 0.4671923527233682
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04618854603468839
$ws.Range("N2").Value = 8.292092835219293
$ws.Range("O2").Value = 4.959649583040046
